# Auto-generated script: applies the 2025-06-05 weekly crime-data refresh
# to output/violent-crime-full-year.xlsx (updates column L = year 2025 YTD totals,
# plus a handful of prior-year reclassification corrections in columns G/K).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 12).Value = 2697   # L2: 2687 -> 2697
$ws.Cells.Item(3, 11).Value = 8184   # K3: 8183 -> 8184
$ws.Cells.Item(3, 12).Value = 2725   # L3: 2712 -> 2725
$ws.Cells.Item(4, 7).Value = 1505   # G4: 1504 -> 1505
$ws.Cells.Item(4, 11).Value = 1765   # K4: 1764 -> 1765
$ws.Cells.Item(4, 12).Value = 731   # L4: 722 -> 731
$ws.Cells.Item(5, 12).Value = 160   # L5: 159 -> 160
$ws.Cells.Item(6, 11).Value = 9121   # K6: 9122 -> 9121
$ws.Cells.Item(6, 12).Value = 2442   # L6: 2424 -> 2442
$ws.Cells.Item(7, 7).Value = 24732   # G7: 24731 -> 24732
$ws.Cells.Item(7, 11).Value = 27557   # K7: 27556 -> 27557
$ws.Cells.Item(7, 12).Value = 8755   # L7: 8704 -> 8755

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(3, 12).Value = 184   # L3: 183 -> 184
$ws.Cells.Item(6, 12).Value = 155   # L6: 154 -> 155
$ws.Cells.Item(7, 12).Value = 556   # L7: 554 -> 556

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 12).Value = 65   # L2: 64 -> 65
$ws.Cells.Item(4, 12).Value = 7   # L4: 6 -> 7
$ws.Cells.Item(7, 12).Value = 207   # L7: 205 -> 207

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(6, 12).Value = 131   # L6: 130 -> 131
$ws.Cells.Item(7, 12).Value = 392   # L7: 391 -> 392

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(6, 12).Value = 24   # L6: 23 -> 24
$ws.Cells.Item(7, 12).Value = 120   # L7: 119 -> 120

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 12).Value = 94   # L3: 93 -> 94
$ws.Cells.Item(4, 12).Value = 20   # L4: 19 -> 20
$ws.Cells.Item(6, 12).Value = 104   # L6: 101 -> 104
$ws.Cells.Item(7, 12).Value = 324   # L7: 319 -> 324

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(3, 12).Value = 7   # L3: 6 -> 7
$ws.Cells.Item(7, 12).Value = 289   # L7: 288 -> 289
$ws.Cells.Item(8, 12).Value = 556   # L8: 554 -> 556
$ws.Cells.Item(13, 12).Value = 13   # L13: 12 -> 13
$ws.Cells.Item(15, 12).Value = 64   # L15: 62 -> 64
$ws.Cells.Item(20, 12).Value = 222   # L20: 221 -> 222
$ws.Cells.Item(21, 12).Value = 24   # L21: 23 -> 24
$ws.Cells.Item(23, 12).Value = 89   # L23: 87 -> 89
$ws.Cells.Item(26, 12).Value = 8   # L26: 7 -> 8
$ws.Cells.Item(29, 12).Value = 464   # L29: 462 -> 464
$ws.Cells.Item(33, 12).Value = 392   # L33: 391 -> 392
$ws.Cells.Item(36, 12).Value = 120   # L36: 118 -> 120
$ws.Cells.Item(37, 12).Value = 324   # L37: 319 -> 324
$ws.Cells.Item(40, 12).Value = 23   # L40: 22 -> 23
$ws.Cells.Item(42, 12).Value = 289   # L42: 287 -> 289
$ws.Cells.Item(43, 12).Value = 68   # L43: 67 -> 68
$ws.Cells.Item(44, 12).Value = 65   # L44: 64 -> 65
$ws.Cells.Item(47, 12).Value = 66   # L47: 65 -> 66
$ws.Cells.Item(49, 12).Value = 47   # L49: 48 -> 47
$ws.Cells.Item(54, 12).Value = 175   # L54: 173 -> 175
$ws.Cells.Item(60, 12).Value = 52   # L60: 51 -> 52
$ws.Cells.Item(63, 7).Value = 306   # G63: 305 -> 306
$ws.Cells.Item(63, 11).Value = 158   # K63: 157 -> 158
$ws.Cells.Item(63, 12).Value = 33   # L63: 29 -> 33
$ws.Cells.Item(67, 12).Value = 323   # L67: 320 -> 323
$ws.Cells.Item(73, 12).Value = 72   # L73: 71 -> 72
$ws.Cells.Item(83, 12).Value = 207   # L83: 205 -> 207
$ws.Cells.Item(84, 12).Value = 92   # L84: 91 -> 92
$ws.Cells.Item(85, 12).Value = 457   # L85: 454 -> 457
$ws.Cells.Item(88, 12).Value = 112   # L88: 111 -> 112
$ws.Cells.Item(90, 12).Value = 85   # L90: 84 -> 85
$ws.Cells.Item(91, 12).Value = 124   # L91: 123 -> 124
$ws.Cells.Item(94, 12).Value = 107   # L94: 106 -> 107
$ws.Cells.Item(95, 12).Value = 120   # L95: 119 -> 120
$ws.Cells.Item(97, 12).Value = 78   # L97: 77 -> 78
$ws.Cells.Item(98, 12).Value = 59   # L98: 58 -> 59
$ws.Cells.Item(101, 7).Value = 24732   # G101: 24731 -> 24732
$ws.Cells.Item(101, 11).Value = 27557   # K101: 27556 -> 27557
$ws.Cells.Item(101, 12).Value = 8755   # L101: 8704 -> 8755

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(4, 12).Value = 26   # L4: 25 -> 26
$ws.Cells.Item(6, 12).Value = 79   # L6: 77 -> 79
$ws.Cells.Item(7, 12).Value = 323   # L7: 320 -> 323

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(2, 12).Value = 34   # L2: 33 -> 34
$ws.Cells.Item(7, 12).Value = 92   # L7: 91 -> 92

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(2, 12).Value = 15   # L2: 16 -> 15
$ws.Cells.Item(7, 12).Value = 47   # L7: 48 -> 47

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 12).Value = 39   # L2: 38 -> 39
$ws.Cells.Item(3, 11).Value = 123   # K3: 122 -> 123
$ws.Cells.Item(6, 11).Value = 278   # K6: 279 -> 278
$ws.Cells.Item(6, 12).Value = 90   # L6: 89 -> 90
$ws.Cells.Item(7, 12).Value = 175   # L7: 173 -> 175

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 12).Value = 172   # L3: 170 -> 172
$ws.Cells.Item(7, 12).Value = 464   # L7: 462 -> 464

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(2, 12).Value = 28   # L2: 27 -> 28
$ws.Cells.Item(7, 12).Value = 65   # L7: 64 -> 65

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 12).Value = 82   # L2: 81 -> 82
$ws.Cells.Item(3, 12).Value = 88   # L3: 87 -> 88
$ws.Cells.Item(7, 12).Value = 289   # L7: 287 -> 289

$ws = $wb.Worksheets.Item('Boystown')
$ws.Cells.Item(5, 12).Value = 6   # L5: 5 -> 6
$ws.Cells.Item(6, 12).Value = 13   # L6: 12 -> 13

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(3, 12).Value = 33   # L3: 32 -> 33
$ws.Cells.Item(6, 12).Value = 18   # L6: 17 -> 18
$ws.Cells.Item(7, 12).Value = 89   # L7: 87 -> 89

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(6, 12).Value = 17   # L6: 16 -> 17
$ws.Cells.Item(7, 12).Value = 124   # L7: 123 -> 124

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Cells.Item(6, 12).Value = 16   # L6: 15 -> 16
$ws.Cells.Item(7, 12).Value = 24   # L7: 23 -> 24

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(4, 12).Value = 20   # L4: 19 -> 20
$ws.Cells.Item(7, 12).Value = 222   # L7: 221 -> 222

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(3, 12).Value = 29   # L3: 28 -> 29
$ws.Cells.Item(6, 12).Value = 30   # L6: 29 -> 30
$ws.Cells.Item(7, 12).Value = 120   # L7: 118 -> 120

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 12).Value = 86   # L2: 85 -> 86
$ws.Cells.Item(7, 12).Value = 289   # L7: 288 -> 289

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(6, 12).Value = 35   # L6: 34 -> 35
$ws.Cells.Item(7, 12).Value = 107   # L7: 106 -> 107

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(3, 12).Value = 22   # L3: 21 -> 22
$ws.Cells.Item(7, 12).Value = 66   # L7: 65 -> 66

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 12).Value = 24   # L2: 23 -> 24
$ws.Cells.Item(5, 12).Value = 1   # L5: None -> 1
$ws.Cells.Item(7, 12).Value = 64   # L7: 62 -> 64

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(4, 12).Value = 6   # L4: 5 -> 6
$ws.Cells.Item(7, 12).Value = 59   # L7: 58 -> 59

$ws = $wb.Worksheets.Item('East Village')
$ws.Cells.Item(4, 12).Value = 2   # L4: 1 -> 2
$ws.Cells.Item(7, 12).Value = 8   # L7: 7 -> 8

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(3, 12).Value = 17   # L3: 16 -> 17
$ws.Cells.Item(7, 12).Value = 72   # L7: 71 -> 72

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(3, 12).Value = 16   # L3: 15 -> 16
$ws.Cells.Item(7, 12).Value = 78   # L7: 77 -> 78

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(6, 12).Value = 34   # L6: 33 -> 34
$ws.Cells.Item(7, 12).Value = 112   # L7: 111 -> 112

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(6, 12).Value = 22   # L6: 21 -> 22
$ws.Cells.Item(7, 12).Value = 85   # L7: 84 -> 85

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(5, 12).Value = 3   # L5: 2 -> 3
$ws.Cells.Item(7, 12).Value = 52   # L7: 51 -> 52

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(2, 12).Value = 15   # L2: 14 -> 15
$ws.Cells.Item(7, 12).Value = 68   # L7: 67 -> 68

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 12).Value = 133   # L2: 132 -> 133
$ws.Cells.Item(3, 12).Value = 185   # L3: 184 -> 185
$ws.Cells.Item(4, 12).Value = 37   # L4: 36 -> 37
$ws.Cells.Item(7, 12).Value = 457   # L7: 454 -> 457

$ws = $wb.Worksheets.Item('Andersonville')
$ws.Cells.Item(2, 12).Value = 2   # L2: 1 -> 2
$ws.Cells.Item(7, 12).Value = 7   # L7: 6 -> 7

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Cells.Item(3, 12).Value = 11   # L3: 10 -> 11
$ws.Cells.Item(7, 12).Value = 23   # L7: 22 -> 23

